$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 25
$ws.Cells.Item($row, 1).Value = 45929
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
$ws.Cells.Item($row, 2).Value = "21,2282"
$ws.Cells.Item($row, 3).Value = "15,0392"
$ws.Cells.Item($row, 4).Value = "15,0392"
$ws.Cells.Item($row, 5).Value = "15,0392"
